# HeroLevelup.xlsx: update the accumulated-experience curve (column D,
# "exp") for levels 2-100 (rows 8-106), and refresh the sheet's current
# selection to reflect the newly-edited range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cumulative-exp values for rows 8..106 (levels 2..100), in order.
$expValues = @(
    12, 28, 52, 87, 137, 241, 381, 563, 794, 1081, 1547, 2107,
    2771, 3548, 4449, 5484, 6663, 7997, 9497, 11174, 13038, 15101, 17373, 19866,
    23273, 26984, 31014, 35378, 40090, 45164, 51706, 58719, 66221, 74230, 82764, 91842,
    101482, 111702, 122521, 133958, 148044, 162895, 178533, 194980, 212258, 232981, 254706, 277458,
    304240, 332237, 364727, 398633, 433988, 470825, 509177, 549077, 590559, 633657, 678405, 724836,
    772984, 822884, 874570, 928076, 983436, 1040686, 1099859, 1160991, 1224116, 1289270, 1356487, 1425803,
    1497253, 1570872, 1646696, 1724760, 1805100, 1887752, 1972751, 2060134, 2149936, 2242194, 2336944, 2434222,
    2534064, 2636507, 2741588, 2849343, 2959809, 3073022, 3189020, 3307839, 3429517, 3554091, 3681597, 3812074,
    3945558, 4082087, 4221698
)

$startRow = 8
for ($i = 0; $i -lt $expValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $expValues[$i]
}

# Match the saved view state: selection now spans D8:D106 with D8 active,
# and the sheet scrolls back to the top (no explicit topLeftCell).
$ws.Range("D8:D106").Select()
